$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$sh.Table.ApplyStyle("{4B9CCE46-CDEB-4F26-9C56-E09C4AADB6A3}")
